$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value = 5149
$ws1.Range("F14").Value = 1316
$ws1.Range("F15").Value = 59
$ws1.Range("F19").Value = 115
$ws1.Range("F21").Value = 178
$ws1.Range("F22").Value = 105
$ws1.Range("F25").Value = 331
$ws1.Range("F27").Value = 3400
$ws1.Range("F29").Value = 2669
$ws1.Range("F31").Value = 1679
$ws1.Range("F32").Value = 3862
$ws1.Range("F34").Value = 909
$ws1.Range("F36").Value = 1218
$ws1.Range("F38").Value = 966
$ws1.Range("F40").Value = 44
$ws1.Range("F41").Value = 963
$ws1.Range("F42").Value = 627
$ws1.Range("F43").Value = 442
$ws1.Range("F44").Value = 385
$ws1.Range("F46").Value = 3528

# Sheet "全部类型" (All types) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 5149
$ws4.Range("F13").Value = 1316
$ws4.Range("F17").Value = 115
$ws4.Range("F20").Value = 178
$ws4.Range("F24").Value = 105
$ws4.Range("F26").Value = 331
$ws4.Range("F27").Value = 3400
$ws4.Range("F31").Value = 2669
$ws4.Range("F32").Value = 1679
$ws4.Range("F33").Value = 3862
$ws4.Range("F36").Value = 909
$ws4.Range("F37").Value = 1218
$ws4.Range("F39").Value = 966
$ws4.Range("F42").Value = 44
$ws4.Range("F43").Value = 963
$ws4.Range("F44").Value = 627
$ws4.Range("F45").Value = 385
$ws4.Range("F49").Value = 3528

